$d = $word.ActiveDocument

$rPr = '<w:rPr><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$pkgOpen  = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body><w:p>'
$pkgClose = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Change 1 (first paragraph): "These Problems could impact the students
# motivation and grades and also causes more work for the professors. "
# -> split into 3 runs with a proofErr-wrapped "students".
# ---------------------------------------------------------------------------
$old1 = "These Problems could impact the students motivation and grades and also causes more work for the professors. "
$f1 = $d.Content
$found1 = $f1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1 = $d.Range($f1.Start, $f1.End)
    $inner1 = '<w:r>' + $rPr + '<w:t xml:space="preserve">These Problems could impact the </w:t></w:r>' +
              '<w:proofErr w:type="gramStart"/>' +
              '<w:r>' + $rPr + '<w:t>students</w:t></w:r>' +
              '<w:proofErr w:type="gramEnd"/>' +
              '<w:r>' + $rPr + '<w:t xml:space="preserve"> motivation and grades and also causes more work for the professors. </w:t></w:r>'
    $r1.InsertXML($pkgOpen + $inner1 + $pkgClose)
}

# ---------------------------------------------------------------------------
# Change 2+3 (final paragraph): the whole paragraph text is rewritten into
# more runs, with a bookmark ("_GoBack") sitting between the two halves.
# Both spans (before/after the bookmark) are replaced together in a single
# InsertXML call so the bookmark markers can be placed explicitly at the
# correct position instead of letting Find/InsertXML relocate them.
# ---------------------------------------------------------------------------
$oldFull = "Our idea of creating a facial recognition software that detects and recognizes students in the lecture room would solve these problems by creating the list at the end of the lecture and marking all students that it recognized as attended. The goal of our project was to create a proof of concept that shows that this idea could be used as an alternative to the attendance list on paper that is currently used. The software could also be improved to save more detailed statistics about the students attendance. "
$f2 = $d.Content
$found2 = $f2.Find.Execute($oldFull, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2 = $d.Range($f2.Start, $f2.End)
    $inner2 = '<w:r>' + $rPr + '<w:t>Our idea of creating a facial recognition software that detects and recognizes students in the lecture room would solve these problems by creating the list at the end of the lecture and marking all students that it recognized as attended.</w:t></w:r>' +
              '<w:r>' + $rPr + '<w:t xml:space="preserve"> An additional advantage of this software would be that we automatically show the names of all the students sitting in the classroom in the camera feed. This helps the professor to remember the names of students and can help to create a more comfortable learning environment for the students.</w:t></w:r>' +
              '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' +
              '<w:proofErr w:type="gramStart"/>' +
              '<w:r>' + $rPr + '<w:t>The</w:t></w:r>' +
              '<w:proofErr w:type="gramEnd"/>' +
              '<w:r>' + $rPr + '<w:t xml:space="preserve"> goal of our project was to create a</w:t></w:r>' +
              '<w:r>' + $rPr + '<w:t xml:space="preserve"> proof of concept that shows that this idea could be used as an alternative to the attendance list on paper that is currently used. The software could </w:t></w:r>' +
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
              '<w:r>' + $rPr + '<w:t xml:space="preserve">also be improved to save more detailed statistics about the </w:t></w:r>' +
              '<w:proofErr w:type="gramStart"/>' +
              '<w:r>' + $rPr + '<w:t>students</w:t></w:r>' +
              '<w:proofErr w:type="gramEnd"/>' +
              '<w:r>' + $rPr + '<w:t xml:space="preserve"> attendance. </w:t></w:r>'
    $r2.InsertXML($pkgOpen + $inner2 + $pkgClose)
}

Write-Output "found1=$found1 found2=$found2"
